$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "juan"
$ws.Range("B4").Value = 42

[void]$ws.Range("B4").Select()
